$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999805566007045
$ws.Range("E2").Value = 0.9999805566007045

# Row 3
$ws.Range("D3").Value = 0.9934485633398517
$ws.Range("E3").Value = 0.9934485633398517

# Row 4
$ws.Range("D4").Value = [double]"7.552765487614542E-38"
$ws.Range("E4").Value = [double]"7.552765487614542E-38"

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.0002091154765998324
$ws.Range("E5").Value = 0.0002091154765998324

# Row 6
$ws.Range("D6").Value = 0.9999999049182763
$ws.Range("E6").Value = 0.9999999049182763

# Row 7
$ws.Range("D7").Value = 0.9999999853018455
$ws.Range("E7").Value = [double]"1.469815447130429E-08"

# Row 8
$ws.Range("D8").Value = 0.9514863751556601
$ws.Range("E8").Value = 0.04851362484433985

# Row 9
$ws.Range("D9").Value = 0.9999999999861076
$ws.Range("E9").Value = [double]"1.389244275173951E-11"

# Row 10
$ws.Range("D10").Value = [double]"3.159487743012221E-07"
$ws.Range("E10").Value = 0.9999996840512257

# Row 11
$ws.Range("D11").Value = [double]"4.510089935594352E-44"
$ws.Range("F11").Value = 14.68696594238281
$ws.Range("G11").Value = 0.5
